$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 4: Root Rush
$ws.Range("H4").Value = 539.8333
$ws.Range("I4").Value = 347.8
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 347.8
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -233.8
$ws.Range("N4").Value = -1728

# Row 5: Met a Sticky End
$ws.Range("H5").Value = 220.22223
$ws.Range("I5").Value = 220.22223
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 220.22223
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -105.22223
$ws.Range("N5").ClearContents()

# Row 17: One for the Road
$ws.Range("H17").Value = 1677.0571
$ws.Range("J17").Value = 1677.0571
$ws.Range("L17").Value = 5031.1713
$ws.Range("N17").Value = -5367.1713

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 209.33333
$ws.Range("I107").Value = 172.26666
$ws.Range("J107").Value = 394.66666
$ws.Range("K107").Value = 172.26666
$ws.Range("L107").Value = 394.66666
$ws.Range("M107").Value = 1747.73334
$ws.Range("N107").Value = -4234.66666

# Row 129: Practical Command
$ws.Range("H129").Value = 1921
$ws.Range("J129").Value = 2335.6667
$ws.Range("L129").Value = 7007.000100000001
$ws.Range("N129").Value = -17007.0001

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 4198.8887
$ws.Range("I141").Value = 3815
$ws.Range("J141").Value = 4966.6665
$ws.Range("K141").Value = 11445
$ws.Range("L141").Value = 14899.9995
$ws.Range("M141").Value = -6265
$ws.Range("N141").Value = -25259.9995

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 1111356.6
$ws.Range("J5").Value = 144
$ws.Range("L5").Value = 144
$ws.Range("N5").Value = -368

# Row 130: A Gift of Gloves
$ws.Range("H130").Value = 88879.5
$ws.Range("J130").Value = 88879.5
$ws.Range("L130").Value = 88879.5
$ws.Range("N130").Value = -98919.5

# Row 131: Additions to the Armoire
$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 4: Mending Fences
$ws.Range("H4").Value = 1111356.6
$ws.Range("J4").Value = 144
$ws.Range("L4").Value = 144
$ws.Range("N4").Value = -374

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2226.0244
$ws.Range("I31").Value = 1685.3125
$ws.Range("J31").Value = 4148.5557
$ws.Range("K31").Value = 1685.3125
$ws.Range("L31").Value = 4148.5557
$ws.Range("M31").Value = -1390.3125
$ws.Range("N31").Value = -4738.5557

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2226.0244
$ws.Range("I34").Value = 1685.3125
$ws.Range("J34").Value = 4148.5557
$ws.Range("K34").Value = 1685.3125
$ws.Range("L34").Value = 4148.5557
$ws.Range("M34").Value = -1483.3125
$ws.Range("N34").Value = -4552.5557

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1578.3334
$ws.Range("I58").Value = 1473.5416
$ws.Range("J58").Value = 2416.6667
$ws.Range("K58").Value = 1473.5416
$ws.Range("L58").Value = 2416.6667
$ws.Range("M58").Value = -1270.5416
$ws.Range("N58").Value = -2822.6667

# Row 75: The Darkest Hearth
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41996

# Row 78: Fruit of the Loom (L)
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129984

# Row 80: The Long Armillae of the Law
$ws.Range("H80").Value = 22500.5
$ws.Range("J80").Value = 22500.5
$ws.Range("L80").Value = 22500.5
$ws.Range("N80").Value = -24746.5

# Row 83: Wooden Ambitions (L)
$ws.Range("H83").Value = 22500.5
$ws.Range("J83").Value = 22500.5
$ws.Range("L83").Value = 67501.5
$ws.Range("N83").Value = -78733.5

# Row 99: O Pine
$ws.Range("H99").Value = 5850332.5
$ws.Range("I99").Value = 8549199
$ws.Range("K99").Value = 8549199
$ws.Range("M99").Value = -8547701

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 2952.8572
$ws.Range("J122").Value = 2892.875
$ws.Range("L122").Value = 8678.625
$ws.Range("N122").Value = -13578.625

# Row 126: A Better Conductor
$ws.Range("H126").Value = 5850332.5
$ws.Range("I126").Value = 8549199
$ws.Range("K126").Value = 25647597
$ws.Range("M126").Value = -25645127

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 35872.793
$ws.Range("I134").Value = 937.6842
$ws.Range("K134").Value = 2813.0526
$ws.Range("M134").Value = -278.0526

# Row 136: Turali Quality
$ws.Range("H136").Value = 1578.3334
$ws.Range("I136").Value = 1473.5416
$ws.Range("J136").Value = 2416.6667
$ws.Range("K136").Value = 4420.6248
$ws.Range("L136").Value = 7250.000100000001
$ws.Range("M136").Value = -1870.6248
$ws.Range("N136").Value = -12350.0001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 103
$ws.Range("J2").Value = 24.25
$ws.Range("L2").Value = 145.5
$ws.Range("N2").Value = -371.5

# Row 5: What a Sap
$ws.Range("H5").Value = 808.7727
$ws.Range("I5").Value = 656.3333
$ws.Range("J5").Value = 991.7
$ws.Range("K5").Value = 1968.9999
$ws.Range("L5").Value = 2975.1
$ws.Range("M5").Value = -1856.9999
$ws.Range("N5").Value = -3199.1

# Row 40: True Grits
$ws.Range("H40").Value = 81.7
$ws.Range("I40").Value = 64.375
$ws.Range("K40").Value = 257.5
$ws.Range("M40").Value = -188.5

# Row 107: Slippery Service
$ws.Range("H107").Value = 328.66666
$ws.Range("I107").Value = 258.33334
$ws.Range("J107").Value = 363.83334
$ws.Range("K107").Value = 775.0000200000001
$ws.Range("L107").Value = 1091.50002
$ws.Range("M107").Value = 1144.99998
$ws.Range("N107").Value = -4931.500019999999

# Row 129: Comfort Food
$ws.Range("H129").Value = 41667370
$ws.Range("J129").Value = 55556350
$ws.Range("L129").Value = 166669050
$ws.Range("N129").Value = -166679050

# Row 130: Blast from the Pasta
$ws.Range("H130").Value = 14980.667
$ws.Range("I130").Value = 17377.8
$ws.Range("K130").Value = 52133.39999999999
$ws.Range("M130").Value = -47113.39999999999

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1413.909
$ws.Range("I131").Value = 1198.1111
$ws.Range("J131").Value = 2385
$ws.Range("K131").Value = 3594.3333
$ws.Range("L131").Value = 7155
$ws.Range("M131").Value = 1445.6667
$ws.Range("N131").Value = -17235

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 808.7727
$ws.Range("I135").Value = 656.3333
$ws.Range("J135").Value = 991.7
$ws.Range("K135").Value = 5906.9997
$ws.Range("L135").Value = 8925.300000000001
$ws.Range("M135").Value = -3371.9997
$ws.Range("N135").Value = -13995.3

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 2023.8125
$ws.Range("I136").Value = 1812.9286
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5438.7858
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -338.7857999999997
$ws.Range("N136").Value = -20700

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 5851.2
$ws.Range("I137").Value = 4468.5713
$ws.Range("J137").Value = 7061
$ws.Range("K137").Value = 13405.7139
$ws.Range("L137").Value = 21183
$ws.Range("M137").Value = -8305.713899999999
$ws.Range("N137").Value = -31383

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 52: It's My Business to Know Things
$ws.Range("H52").Value = 19993.334

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 8433.333000000001
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 11400
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 34200
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -39140

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 6: Sticking Their Necks Out
$ws.Range("H6").Value = 75000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Row 40: Best Served Toad
$ws.Range("H40").Value = 3972309
$ws.Range("I40").Value = 4282.273
$ws.Range("J40").Value = 18521740
$ws.Range("K40").Value = 4282.273
$ws.Range("L40").Value = 18521740
$ws.Range("M40").Value = -4146.273
$ws.Range("N40").Value = -18522012

# Row 120: Into the Storm
$ws.Range("H120").Value = 20698
$ws.Range("J120").Value = 20698
$ws.Range("L120").Value = 20698
$ws.Range("N120").Value = -30374

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 29: Getting Handsy
$ws.Range("H29").Value = 275000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 9810.556
$ws.Range("I62").Value = 9299.333000000001
$ws.Range("J62").Value = 10833
$ws.Range("K62").Value = 9299.333000000001
$ws.Range("L62").Value = 10833
$ws.Range("M62").Value = -8675.333000000001
$ws.Range("N62").Value = -12081

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 9810.556
$ws.Range("I65").Value = 9299.333000000001
$ws.Range("J65").Value = 10833
$ws.Range("K65").Value = 46496.665
$ws.Range("L65").Value = 54165
$ws.Range("M65").Value = -43376.665
$ws.Range("N65").Value = -60405

# Row 100: Of Great Import
$ws.Range("H100").Value = 4763484
$ws.Range("I100").Value = 5496213.5
$ws.Range("J100").Value = 741.5
$ws.Range("K100").Value = 10992427
$ws.Range("L100").Value = 1483
$ws.Range("M100").Value = -10991886
$ws.Range("N100").Value = -2565

# Row 101: Who War It Better
$ws.Range("H101").Value = 63301
$ws.Range("J101").Value = 63301
$ws.Range("L101").Value = 63301
$ws.Range("N101").Value = -69791

# Row 110: Suits You
$ws.Range("H110").Value = 60097.4
$ws.Range("J110").Value = 60097.4
$ws.Range("L110").Value = 60097.4
$ws.Range("N110").Value = -68277.39999999999

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 4780.6
$ws.Range("I126").Value = 3702
$ws.Range("J126").Value = 5499.6665
$ws.Range("K126").Value = 11106
$ws.Range("L126").Value = 16498.9995
$ws.Range("M126").Value = -8636
$ws.Range("N126").Value = -21438.9995

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2155.9688
$ws.Range("I132").Value = 1831.2273
$ws.Range("J132").Value = 2870.4
$ws.Range("K132").Value = 5493.6819
$ws.Range("L132").Value = 8611.200000000001
$ws.Range("M132").Value = -2963.6819
$ws.Range("N132").Value = -13671.2
